# Regenerate save_data to use K (strikeouts) instead of Strike# values.
# Updates column G ("K") for rows 2-24 on the active sheet with the
# recalculated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 4
    4  = 3
    5  = 6
    6  = 3
    7  = 0
    8  = 3
    9  = 4
    10 = 4
    11 = 4
    12 = 5
    13 = 0
    14 = 0
    15 = 2
    16 = 2
    17 = 9
    18 = 2
    19 = 6
    20 = 2
    21 = 3
    22 = 4
    23 = 2
    24 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
